# Auto-generated: apply scheduled-runner market-data refresh to Leve profit sheets.
# For each changed row, updates currentAveragePrice / *NQ / *HQ, LevePriceNQ/HQ,
# and LeveProfitNQ/HQ (columns H-N) to the latest computed figures. Cells that the
# refresh leaves blank are cleared; cells that newly have a computed profit are set.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Range("H40").Value = 3529.4546
$ws.Range("I40").Value = 2820
$ws.Range("J40").Value = 3934.8572
$ws.Range("K40").Value = 2820
$ws.Range("L40").Value = 3934.8572
$ws.Range("M40").Value = -2645
$ws.Range("N40").Value = -4284.8572
# Row 51
$ws.Range("H51").Value = 9867
$ws.Range("J51").Value = 9974
$ws.Range("L51").Value = 9974
$ws.Range("N51").Value = -10942
# Row 64
$ws.Range("H64").Value = 9912
$ws.Range("I64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("M64").Value = $null
# Row 67
$ws.Range("H67").Value = 9912
$ws.Range("I67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("M67").Value = $null
# Row 116
$ws.Range("H116").Value = 5860.3335
$ws.Range("I116").Value = 6159.6
$ws.Range("J116").Value = 5646.5713
$ws.Range("K116").Value = 6159.6
$ws.Range("L116").Value = 5646.5713
$ws.Range("M116").Value = -2717.6
$ws.Range("N116").Value = -12530.5713

$ws = $wb.Worksheets.Item("ARM")
# Row 5
$ws.Range("H5").Value = 236
$ws.Range("I5").Value = 236
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 236
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -124
$ws.Range("N5").Value = $null
# Row 97
$ws.Range("H97").Value = 1209.5333
$ws.Range("I97").Value = 395.69232
$ws.Range("K97").Value = 395.69232
$ws.Range("M97").Value = 100.30768
# Row 132
$ws.Range("H132").Value = 3050.111
$ws.Range("I132").Value = 2372
$ws.Range("K132").Value = 7116
$ws.Range("M132").Value = -4586
# Row 135
$ws.Range("H135").Value = 88699.5
$ws.Range("J135").Value = 88699.5
$ws.Range("L135").Value = 88699.5
$ws.Range("N135").Value = -98839.5

$ws = $wb.Worksheets.Item("BSM")
# Row 4
$ws.Range("H4").Value = 236
$ws.Range("I4").Value = 236
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 236
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -121
$ws.Range("N4").Value = $null
# Row 10
$ws.Range("H10").Value = 1502.5
$ws.Range("I10").Value = 1005
$ws.Range("J10").Value = 2000
$ws.Range("K10").Value = 1005
$ws.Range("L10").Value = 2000
$ws.Range("M10").Value = -865
$ws.Range("N10").Value = -2280
# Row 94
$ws.Range("H94").Value = 2466.2
$ws.Range("I94").Value = 1133.625
$ws.Range("J94").Value = 3989.1428
$ws.Range("K94").Value = 1133.625
$ws.Range("L94").Value = 3989.1428
$ws.Range("M94").Value = -682.625
$ws.Range("N94").Value = -4891.1428
# Row 105
$ws.Range("H105").Value = 3865.7778
$ws.Range("I105").Value = 3905
$ws.Range("J105").Value = 3854.5715
$ws.Range("K105").Value = 3905
$ws.Range("L105").Value = 3854.5715
$ws.Range("M105").Value = -2158
$ws.Range("N105").Value = -7348.5715
# Row 135
$ws.Range("H135").Value = 99999
$ws.Range("J135").Value = 99999
$ws.Range("L135").Value = 99999
$ws.Range("N135").Value = -110139

$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 660
$ws.Range("I22").Value = 800
$ws.Range("K22").Value = 800
$ws.Range("M22").Value = -450
# Row 31
$ws.Range("H31").Value = 3314.7
$ws.Range("I31").Value = 2013
$ws.Range("K31").Value = 2013
$ws.Range("M31").Value = -1718
# Row 34
$ws.Range("H34").Value = 3314.7
$ws.Range("I34").Value = 2013
$ws.Range("K34").Value = 2013
$ws.Range("M34").Value = -1811
# Row 88
$ws.Range("H88").Value = 35781
$ws.Range("J88").Value = 35781
$ws.Range("L88").Value = 35781
$ws.Range("N88").Value = -36593
# Row 91
$ws.Range("H91").Value = 35781
$ws.Range("J91").Value = 35781
$ws.Range("L91").Value = 35781
$ws.Range("N91").Value = -38589

$ws = $wb.Worksheets.Item("CUL")
# Row 23
$ws.Range("H23").Value = 428.83334
$ws.Range("J23").Value = 469.25
$ws.Range("L23").Value = 1407.75
$ws.Range("N23").Value = -1877.75
# Row 38
$ws.Range("H38").Value = 138.88889
$ws.Range("I38").Value = 37.75
$ws.Range("J38").Value = 219.8
$ws.Range("K38").Value = 113.25
$ws.Range("L38").Value = 659.4000000000001
$ws.Range("M38").Value = 233.75
$ws.Range("N38").Value = -1353.4
# Row 92
$ws.Range("H92").Value = 313.44446
$ws.Range("I92").Value = 550
$ws.Range("K92").Value = 1650
$ws.Range("M92").Value = -402
# Row 98
$ws.Range("H98").Value = 417
$ws.Range("J98").Value = 344.125
$ws.Range("L98").Value = 1032.375
$ws.Range("N98").Value = -4028.375

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 12360.056
$ws.Range("I80").Value = 4828.375
$ws.Range("J80").Value = 18385.4
$ws.Range("K80").Value = 4828.375
$ws.Range("L80").Value = 18385.4
$ws.Range("M80").Value = -3830.375
$ws.Range("N80").Value = -20381.4
# Row 83
$ws.Range("H83").Value = 12360.056
$ws.Range("I83").Value = 4828.375
$ws.Range("J83").Value = 18385.4
$ws.Range("K83").Value = 24141.875
$ws.Range("L83").Value = 91927
$ws.Range("M83").Value = -19149.875
$ws.Range("N83").Value = -101911
# Row 122
$ws.Range("H122").Value = 1453.5
$ws.Range("I122").Value = 1453.5
$ws.Range("K122").Value = 4360.5
$ws.Range("M122").Value = -1910.5
# Row 126
$ws.Range("H126").Value = 2967.0833
$ws.Range("I126").Value = 2944.8333
$ws.Range("K126").Value = 8834.499899999999
$ws.Range("M126").Value = -6364.499899999999

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 2483.6924
$ws.Range("I7").Value = 1414.2222
$ws.Range("K7").Value = 1414.2222
$ws.Range("M7").Value = -1302.2222
# Row 22
$ws.Range("H22").Value = 870.75
$ws.Range("I22").Value = 870.75
$ws.Range("K22").Value = 870.75
$ws.Range("M22").Value = -575.75
# Row 27
$ws.Range("H27").Value = 870.75
$ws.Range("I27").Value = 870.75
$ws.Range("K27").Value = 870.75
$ws.Range("M27").Value = -763.75
# Row 40
$ws.Range("H40").Value = 3040.8
$ws.Range("I40").Value = 2840.5386
$ws.Range("K40").Value = 2840.5386
$ws.Range("M40").Value = -2704.5386
# Row 68
$ws.Range("H68").Value = 2535
$ws.Range("I68").Value = 2672
$ws.Range("J68").Value = 2032.6666
$ws.Range("K68").Value = 2672
$ws.Range("L68").Value = 2032.6666
$ws.Range("M68").Value = -1923
$ws.Range("N68").Value = -3530.6666
# Row 71
$ws.Range("H71").Value = 2535
$ws.Range("I71").Value = 2672
$ws.Range("J71").Value = 2032.6666
$ws.Range("K71").Value = 13360
$ws.Range("L71").Value = 10163.333
$ws.Range("M71").Value = -9616
$ws.Range("N71").Value = -17651.333
# Row 107
$ws.Range("H107").Value = 1039.5
$ws.Range("I107").Value = 1039.5
$ws.Range("K107").Value = 1039.5
$ws.Range("M107").Value = 880.5
# Row 126
$ws.Range("H126").Value = 2483.6924
$ws.Range("I126").Value = 1414.2222
$ws.Range("K126").Value = 4242.6666
$ws.Range("M126").Value = -1772.6666

$ws = $wb.Worksheets.Item("WVR")
# Row 126
$ws.Range("H126").Value = 8499.666999999999
$ws.Range("I126").Value = 12000
$ws.Range("J126").Value = 1499
$ws.Range("K126").Value = 36000
$ws.Range("L126").Value = 4497
$ws.Range("M126").Value = -33530
$ws.Range("N126").Value = -9437
